$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(40).Insert()

$ws.Cells.Item(40, 1).Value = 8
$ws.Cells.Item(40, 2).Value = 'Terminal La Palmera de La Serena'
$ws.Cells.Item(40, 3).Value = 'Coquimbo'
$ws.Cells.Item(40, 4).Value = 44498
$ws.Cells.Item(40, 5).Value = 4
$ws.Cells.Item(40, 6).Value = 100112031
$ws.Cells.Item(40, 7).Value = 'Poroto verde'
$ws.Cells.Item(40, 8).Value = 'Magnum'
$ws.Cells.Item(40, 9).Value = 'Primera'
$ws.Cells.Item(40, 10).Value = 500
$ws.Cells.Item(40, 11).Value = 34500
$ws.Cells.Item(40, 12).Value = 35000
$ws.Cells.Item(40, 13).Value = 34750
$ws.Cells.Item(40, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(40, 15).Value = 'Perú'
$ws.Cells.Item(40, 16).Value = 1390
$ws.Cells.Item(40, 17).Value = 25
$ws.Cells.Item(40, 18).Value = 'Hortaliza'
